$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new header cells L1:T1 ---
$ws.Range("L1").Value = "issueDate"
$ws.Range("M1").Value = "AFM_issuer"
$ws.Range("N1").Value = "AA"
$ws.Range("O1").Value = "series"
$ws.Range("P1").Value = "totalNetValue"
$ws.Range("Q1").Value = "totalVatAmount"
$ws.Range("R1").Value = "totalValue"
$ws.Range("S1").Value = "character"
$ws.Range("T1").Value = "type"

# copy the existing header style (bold, bordered, centered) onto the new header cells
$ws.Range("K1").Copy()
$ws.Range("L1:T1").PasteSpecial(-4122)

# --- Row 2: new trailing empty cells L2:T2 ---
$ws.Range("L2:T2").Style = "Normal"

# --- Row 3: brand-new data row ---
$ws.Range("A3").Value = "'400011184530011"
$ws.Range("B3").Value = "'094439854"
$ws.Range("C3").Value = "'ΤΡΑΚΑΔΑΣ Α.Ε."
$ws.Range("D3").Value = "'5Μ0ΤΔΑ"
$ws.Range("E3").Value = "'4626"
$ws.Range("F3").Value = "'04/10/2025"
$ws.Range("G3").Value = "'1.1"
$ws.Range("I3").Value = "'45,54"
$ws.Range("J3").Value = "'10,93"
$ws.Range("K3").Value = "'56,47"

# H3 stays empty but present, as do L3:T3
$ws.Range("A3:T3").Style = "Normal"

Write-Host "done"
